$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rework sheet "PID4CatRecord" -> "HandleAPIRecord"
# ---------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("PID4CatRecord")
$sheet1.Cells.Validation.Delete()
$sheet1.Cells.Clear()
$sheet1.Name = "HandleAPIRecord"
$sheet1.Range("A1").Value = "response_code"
$sheet1.Range("B1").Value = "handle"
$sheet1.Range("C1").Value = "values"

# ---------------------------------------------------------------
# 2. Insert new sheet "HandleRecord" right after HandleAPIRecord
# ---------------------------------------------------------------
$handleRecord = $wb.Worksheets.Add($null, $sheet1)
$handleRecord.Name = "HandleRecord"
$handleRecord.Range("A1").Value = "index"
$handleRecord.Range("B1").Value = "type"
$handleRecord.Range("C1").Value = "data"
$handleRecord.Range("D1").Value = "ttl"
$handleRecord.Range("E1").Value = "timestamp"
$handleRecord.Range("B2:B1048576").Validation.Add(3, 1, 1, '"URL,STATUS,SCHEMA_VER,LICENSE,EMAIL,RESOURCE_INFO,RELATED,LOG"')

# ---------------------------------------------------------------
# 3. Insert new sheet "HandleData" right after HandleRecord
# ---------------------------------------------------------------
$handleData = $wb.Worksheets.Add($null, $handleRecord)
$handleData.Name = "HandleData"
$handleData.Range("A1").Value = "format"
$handleData.Range("B1").Value = "value"

# ---------------------------------------------------------------
# 4. Move "Container" right after HandleData, rename to HandleRecordContainer
# ---------------------------------------------------------------
$container = $wb.Worksheets.Item("Container")
$container.Move($null, $handleData)
# NOTE: re-fetch by name; the COM reference resolves by position and
# becomes stale (pointing at whatever sheet now occupies the old slot)
# immediately after a Move operation.
$container = $wb.Worksheets.Item("Container")
$container.Name = "HandleRecordContainer"

# ---------------------------------------------------------------
# Remaining sheets (PID4CatRelation, ResourceInfo, LogRecord, Agent,
# RepresentationVariant) are left untouched; they keep their content
# and naturally shift position because of the inserts/moves above.
# ---------------------------------------------------------------
